$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "State Enrolled" (column G) from "N" to "Y" for rows 6, 7, 8, and 10
$ws.Range("G6").Value = "Y"
$ws.Range("G7").Value = "Y"
$ws.Range("G8").Value = "Y"
$ws.Range("G10").Value = "Y"
